# Update rules in DiscountRules.xlsx
#
# The decision table on the (only) sheet currently has 3 columns:
#   A = NAME, B = CONDITION, C = ACTION
# This change inserts two more CONDITION columns between the existing
# CONDITION column and the ACTION column (so ACTION moves from C to E),
# fills in the new column headers/template row, and appends a brand new
# rule row ("Rule 3") that only populates the 3rd condition column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch D1:E26 so the two new columns get persisted (blank) cells for
# every existing row, matching the shift caused by inserting two new
# condition columns before the old ACTION column.
$ws.Range("D1:E26").Style = "Normal"

# The old ACTION column data (currently in C, rows 18 & 19 are the only
# non-blank cells in that column) moves two columns to the right, into
# E, to make room for the two new CONDITION columns.
$ws.Range("E18").Value2 = $ws.Range("C18").Value2
$ws.Range("E19").Value2 = $ws.Range("C19").Value2

# Row 18 holds the column-type headers: NAME | CONDITION | CONDITION | CONDITION | ACTION
$ws.Range("C18").Value = "CONDITION"
$ws.Range("D18").Value = "CONDITION"

# Row 19 holds the per-column template/sample values.
$ws.Range("C19").Value = "s"
$ws.Range("D19").Value = "3rd condition"

# New rule row 27: only the name and the 3rd condition column are
# populated; force the remaining cells of the row to persist as blank,
# matching the style of the other rule rows.
$ws.Range("A27:E27").Style = "Normal"
$ws.Range("A27").Value = "Rule 3"
$ws.Range("D27").Value = "3rd Rule"
